# Weekly fruit/vegetable price update: insert two new daily records at the
# top of the data block (row 26) and push the existing history down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 26 -- this shifts the existing rows 26:112
# down to 28:114 (and the sheet dimension grows to A1:R114 automatically).
$ws.Rows.Item(26).EntireRow.Insert()
$ws.Rows.Item(26).EntireRow.Insert()

# Populate the newly inserted row 26 with the new weekly record.
$ws.Cells.Item(26, 1).Value  = 9
$ws.Cells.Item(26, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(26, 3).Value  = "Metropolitana"
$ws.Cells.Item(26, 4).Value  = 44487
$ws.Cells.Item(26, 5).Value  = 13
$ws.Cells.Item(26, 6).Value  = 100112003
$ws.Cells.Item(26, 7).Value  = "Ajo"
$ws.Cells.Item(26, 8).Value  = "Chino"
$ws.Cells.Item(26, 9).Value  = "Primera"
$ws.Cells.Item(26, 10).Value = 520
$ws.Cells.Item(26, 11).Value = 17000
$ws.Cells.Item(26, 12).Value = 18000
$ws.Cells.Item(26, 13).Value = 17500
$ws.Cells.Item(26, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(26, 15).Value = "China"
$ws.Cells.Item(26, 16).Value = 1750
$ws.Cells.Item(26, 17).Value = 10
$ws.Cells.Item(26, 18).Value = "Hortaliza"

# Populate the newly inserted row 27 with the new weekly record.
$ws.Cells.Item(27, 1).Value  = 9
$ws.Cells.Item(27, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(27, 3).Value  = "Metropolitana"
$ws.Cells.Item(27, 4).Value  = 44487
$ws.Cells.Item(27, 5).Value  = 13
$ws.Cells.Item(27, 6).Value  = 100112003
$ws.Cells.Item(27, 7).Value  = "Ajo"
$ws.Cells.Item(27, 8).Value  = "Chino"
$ws.Cells.Item(27, 9).Value  = "Primera"
$ws.Cells.Item(27, 10).Value = 340
$ws.Cells.Item(27, 11).Value = 18000
$ws.Cells.Item(27, 12).Value = 19000
$ws.Cells.Item(27, 13).Value = 18500
$ws.Cells.Item(27, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(27, 15).Value = "China"
$ws.Cells.Item(27, 16).Value = 1850
$ws.Cells.Item(27, 17).Value = 10
$ws.Cells.Item(27, 18).Value = "Hortaliza"

Write-Host "Inserted weekly Ajo price rows at 26-27; sheet now spans to row 114."
